$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.202.03'
$ws.Range('E2').Value = '  -1.55%  '
$ws.Range('D3').Value = '1.660.49'
$ws.Range('E3').Value = '  -1.55%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.28'
$ws.Range('E5').Value = '  +0.24%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5213'
$ws.Range('E6').Value = '  -2.18%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.005'
$ws.Range('E7').Value = '  +0.29%  '
$ws.Range('E8').Value = '  -0.72%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06320'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.09'
$ws.Range('E10').Value = '  -2.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07724'
$ws.Range('E11').Value = '  -0.98%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.661.57'
$ws.Range('E12').Value = '  -1.62%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.429'
$ws.Range('E13').Value = '  -1.83%  '
$ws.Range('D14').Value = '1.884.79'
$ws.Range('E14').Value = '  -1.72%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5465'
$ws.Range('E15').Value = '  -3.26%  '
$ws.Range('D16').Value = '0.0₅8228'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '64.95'
$ws.Range('E17').Value = '  -2.25%  '
$ws.Range('D18').Value = '26.232.02'
$ws.Range('E18').Value = '  -1.61%  '
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.662'
$ws.Range('E20').Value = '  -3.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.21'
$ws.Range('E21').Value = '  -1.42%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.15'
$ws.Range('E22').Value = '  -2.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.096'
$ws.Range('E23').Value = '  -4.62%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.007'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '138.28'
$ws.Range('E25').Value = '  -3.77%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1241'
$ws.Range('E26').Value = '  -3.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.214'
$ws.Range('E27').Value = '  -3.63%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '16.19'
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.420'
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.06005'
$ws.Range('E30').Value = '  -3.03%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.284'
$ws.Range('E31').Value = '  +0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.571'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.331'
$ws.Range('E33').Value = '  -4.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.645'
$ws.Range('E34').Value = '  -3.46%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9813'
$ws.Range('E35').Value = '  -3.35%  '
$ws.Range('B36').Value = 'MXToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.781'
$ws.Range('E36').Value = '  -0.50%  '
$ws.Range('B37').Value = 'HuobiToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.407'
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5881'
$ws.Range('E38').Value = '  +2.80%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01585'
$ws.Range('E39').Value = '  -3.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.941'
$ws.Range('E40').Value = '  -1.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8633'
$ws.Range('E41').Value = '  -0.36%  '
$ws.Range('E42').Value = '  +0.23%  '
$ws.Range('D43').Value = '1.038.64'
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '99.46'
$ws.Range('E44').Value = '  -1.06%  '
$ws.Range('D45').Value = '1.798.42'
$ws.Range('E45').Value = '  -2.15%  '
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').Value = '0.0₈110'
$ws.Range('E46').Value = '  +0.55%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '57.06'
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.006'
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.064'
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('E50').Value = '  -0.73%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4230'
$ws.Range('E51').Value = '  -0.24%  '
